$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows (50 and 51) by copying formatting from the last existing data row (49)
$ws.Range("A49:E49").Copy()
$ws.Range("A50:E51").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update all data rows (2-51) with refreshed confidence score results
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 0.19109
$ws.Cells.Item(2,3).Value = 0.19109
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(3,1).Value = 16
$ws.Cells.Item(3,2).Value = 0.15315
$ws.Cells.Item(3,3).Value = 0.15315
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(4,1).Value = 44
$ws.Cells.Item(4,2).Value = 0.15315
$ws.Cells.Item(4,3).Value = 0.15315
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(5,1).Value = 39
$ws.Cells.Item(5,2).Value = 0.15315
$ws.Cells.Item(5,3).Value = 0.15315
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(6,1).Value = 35
$ws.Cells.Item(6,2).Value = 0.22368
$ws.Cells.Item(6,3).Value = 0.2236799999999999
$ws.Cells.Item(6,4).Value = [double]"8.326672684688674e-17"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(7,1).Value = 24
$ws.Cells.Item(7,2).Value = 0.18266
$ws.Cells.Item(7,3).Value = 0.18266
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(8,1).Value = 20
$ws.Cells.Item(8,2).Value = 0.15315
$ws.Cells.Item(8,3).Value = 0.15315
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(9,1).Value = 25
$ws.Cells.Item(9,2).Value = 0.15315
$ws.Cells.Item(9,3).Value = 0.15315
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(10,1).Value = 5
$ws.Cells.Item(10,2).Value = 0.22762
$ws.Cells.Item(10,3).Value = 0.2275496
$ws.Cells.Item(10,4).Value = [double]"7.0400000000026e-05"
$ws.Cells.Item(10,5).Value = 0.9999296049558111
$ws.Cells.Item(11,1).Value = 26
$ws.Cells.Item(11,2).Value = 0.22481
$ws.Cells.Item(11,3).Value = 0.224928
$ws.Cells.Item(11,4).Value = 0.0001179999999999792
$ws.Cells.Item(11,5).Value = 0.9998820139223571
$ws.Cells.Item(12,1).Value = 15
$ws.Cells.Item(12,2).Value = 0.16861
$ws.Cells.Item(12,3).Value = 0.1687672666666667
$ws.Cells.Item(12,4).Value = 0.0001572666666666833
$ws.Cells.Item(12,5).Value = 0.9998427580622486
$ws.Cells.Item(13,1).Value = 36
$ws.Cells.Item(13,2).Value = 0.20711
$ws.Cells.Item(13,3).Value = 0.2072984
$ws.Cells.Item(13,4).Value = 0.000188400000000033
$ws.Cells.Item(13,5).Value = 0.999811635487874
$ws.Cells.Item(14,1).Value = 32
$ws.Cells.Item(14,2).Value = 0.16186
$ws.Cells.Item(14,3).Value = 0.1616671583333334
$ws.Cells.Item(14,4).Value = 0.0001928416666666377
$ws.Cells.Item(14,5).Value = 0.9998071955140718
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 0.17142
$ws.Cells.Item(15,3).Value = 0.1716557209523809
$ws.Cells.Item(15,4).Value = 0.0002357209523809467
$ws.Cells.Item(15,5).Value = 0.9997643345988919
$ws.Cells.Item(16,1).Value = 6
$ws.Cells.Item(16,2).Value = 0.21638
$ws.Cells.Item(16,3).Value = 0.2167172
$ws.Cells.Item(16,4).Value = 0.0003372000000000375
$ws.Cells.Item(16,5).Value = 0.9996629136655119
$ws.Cells.Item(17,1).Value = 23
$ws.Cells.Item(17,2).Value = 0.17732
$ws.Cells.Item(17,3).Value = 0.1769598285714286
$ws.Cells.Item(17,4).Value = 0.0003601714285714486
$ws.Cells.Item(17,5).Value = 0.9996399582481806
$ws.Cells.Item(18,1).Value = 10
$ws.Cells.Item(18,2).Value = 0.16776
$ws.Cells.Item(18,3).Value = 0.1682296
$ws.Cells.Item(18,4).Value = 0.0004695999999999867
$ws.Cells.Item(18,5).Value = 0.9995306204206506
$ws.Cells.Item(19,1).Value = 18
$ws.Cells.Item(19,2).Value = 0.16776
$ws.Cells.Item(19,3).Value = 0.1682296
$ws.Cells.Item(19,4).Value = 0.0004695999999999867
$ws.Cells.Item(19,5).Value = 0.9995306204206506
$ws.Cells.Item(20,1).Value = 43
$ws.Cells.Item(20,2).Value = 0.20682
$ws.Cells.Item(20,3).Value = 0.2072984
$ws.Cells.Item(20,4).Value = 0.0004784000000000177
$ws.Cells.Item(20,5).Value = 0.9995218287571226
$ws.Cells.Item(21,1).Value = 8
$ws.Cells.Item(21,2).Value = 0.16158
$ws.Cells.Item(21,3).Value = 0.1610218333333333
$ws.Cells.Item(21,4).Value = 0.0005581666666666929
$ws.Cells.Item(21,5).Value = 0.9994421447095612
$ws.Cells.Item(22,1).Value = 41
$ws.Cells.Item(22,2).Value = 0.18462
$ws.Cells.Item(22,3).Value = 0.1839909
$ws.Cells.Item(22,4).Value = 0.0006290999999999936
$ws.Cells.Item(22,5).Value = 0.9993712955179895
$ws.Cells.Item(23,1).Value = 28
$ws.Cells.Item(23,2).Value = 0.16018
$ws.Cells.Item(23,3).Value = 0.15952486
$ws.Cells.Item(23,4).Value = 0.0006551399999999985
$ws.Cells.Item(23,5).Value = 0.999345288927412
$ws.Cells.Item(24,1).Value = 12
$ws.Cells.Item(24,2).Value = 0.15877
$ws.Cells.Item(24,3).Value = 0.1595909
$ws.Cells.Item(24,4).Value = 0.000820900000000041
$ws.Cells.Item(24,5).Value = 0.9991797733240781
$ws.Cells.Item(25,1).Value = 19
$ws.Cells.Item(25,2).Value = 0.17732
$ws.Cells.Item(25,3).Value = 0.1764531
$ws.Cells.Item(25,4).Value = 0.0008669000000000038
$ws.Cells.Item(25,5).Value = 0.9991338508646853
$ws.Cells.Item(26,1).Value = 42
$ws.Cells.Item(26,2).Value = 0.17535
$ws.Cells.Item(26,3).Value = 0.1744585833333333
$ws.Cells.Item(26,4).Value = 0.0008914166666666723
$ws.Cells.Item(26,5).Value = 0.9991093772492969
$ws.Cells.Item(27,1).Value = 7
$ws.Cells.Item(27,2).Value = 0.17704
$ws.Cells.Item(27,3).Value = 0.1761309114285714
$ws.Cells.Item(27,4).Value = 0.0009090885714285735
$ws.Cells.Item(27,5).Value = 0.9990917371199757
$ws.Cells.Item(28,1).Value = 14
$ws.Cells.Item(28,2).Value = 0.16158
$ws.Cells.Item(28,3).Value = 0.1606594777777778
$ws.Cells.Item(28,4).Value = 0.0009205222222222476
$ws.Cells.Item(28,5).Value = 0.999080324359642
$ws.Cells.Item(29,1).Value = 9
$ws.Cells.Item(29,2).Value = 0.17535
$ws.Cells.Item(29,3).Value = 0.1743629833333333
$ws.Cells.Item(29,4).Value = 0.0009870166666666735
$ws.Cells.Item(29,5).Value = 0.9990139565746283
$ws.Cells.Item(30,1).Value = 17
$ws.Cells.Item(30,2).Value = 0.15456
$ws.Cells.Item(30,3).Value = 0.15351321
$ws.Cells.Item(30,4).Value = 0.00104679000000002
$ws.Cells.Item(30,5).Value = 0.9989543046234632
$ws.Cells.Item(31,1).Value = 22
$ws.Cells.Item(31,2).Value = 0.1672
$ws.Cells.Item(31,3).Value = 0.1660655
$ws.Cells.Item(31,4).Value = 0.001134500000000011
$ws.Cells.Item(31,5).Value = 0.9988667856317008
$ws.Cells.Item(32,1).Value = 40
$ws.Cells.Item(32,2).Value = 0.18378
$ws.Cells.Item(32,3).Value = 0.1824759
$ws.Cells.Item(32,4).Value = 0.001304100000000002
$ws.Cells.Item(32,5).Value = 0.9986975984618459
$ws.Cells.Item(33,1).Value = 29
$ws.Cells.Item(33,2).Value = 0.15737
$ws.Cells.Item(33,3).Value = 0.1592726833333333
$ws.Cells.Item(33,4).Value = 0.001902683333333322
$ws.Cells.Item(33,5).Value = 0.998100929995513
$ws.Cells.Item(34,1).Value = 46
$ws.Cells.Item(34,2).Value = 0.20345
$ws.Cells.Item(34,3).Value = 0.2063149
$ws.Cells.Item(34,4).Value = 0.002864900000000031
$ws.Cells.Item(34,5).Value = 0.9971432842050808
$ws.Cells.Item(35,1).Value = 30
$ws.Cells.Item(35,2).Value = 0.1939
$ws.Cells.Item(35,3).Value = 0.1973001
$ws.Cells.Item(35,4).Value = 0.003400099999999989
$ws.Cells.Item(35,5).Value = 0.9966114215057384
$ws.Cells.Item(36,1).Value = 4
$ws.Cells.Item(36,2).Value = 0.16439
$ws.Cells.Item(36,3).Value = 0.1605872837301587
$ws.Cells.Item(36,4).Value = 0.0038027162698413
$ws.Cells.Item(36,5).Value = 0.9962116895997529
$ws.Cells.Item(37,1).Value = 21
$ws.Cells.Item(37,2).Value = 0.16861
$ws.Cells.Item(37,3).Value = 0.1733706377627928
$ws.Cells.Item(37,4).Value = 0.004760637762792758
$ws.Cells.Item(37,5).Value = 0.9952619185267919
$ws.Cells.Item(38,1).Value = 11
$ws.Cells.Item(38,2).Value = 0.17985
$ws.Cells.Item(38,3).Value = 0.1846128
$ws.Cells.Item(38,4).Value = 0.004762800000000011
$ws.Cells.Item(38,5).Value = 0.9952597767353648
$ws.Cells.Item(39,1).Value = 49
$ws.Cells.Item(39,2).Value = 0.21244
$ws.Cells.Item(39,3).Value = 0.21727
$ws.Cells.Item(39,4).Value = 0.004830000000000029
$ws.Cells.Item(39,5).Value = 0.9951932167630344
$ws.Cells.Item(40,1).Value = 33
$ws.Cells.Item(40,2).Value = 0.16748
$ws.Cells.Item(40,3).Value = 0.1623212
$ws.Cells.Item(40,4).Value = 0.005158799999999991
$ws.Cells.Item(40,5).Value = 0.9948676766298021
$ws.Cells.Item(41,1).Value = 48
$ws.Cells.Item(41,2).Value = 0.17142
$ws.Cells.Item(41,3).Value = 0.1658748133333333
$ws.Cells.Item(41,4).Value = 0.005545186666666646
$ws.Cells.Item(41,5).Value = 0.9944853928593217
$ws.Cells.Item(42,1).Value = 38
$ws.Cells.Item(42,2).Value = 0.18153
$ws.Cells.Item(42,3).Value = 0.1757457333333333
$ws.Cells.Item(42,4).Value = 0.005784266666666676
$ws.Cells.Item(42,5).Value = 0.9942489986586918
$ws.Cells.Item(43,1).Value = 47
$ws.Cells.Item(43,2).Value = 0.18153
$ws.Cells.Item(43,3).Value = 0.1752486095238095
$ws.Cells.Item(43,4).Value = 0.006281390476190474
$ws.Cells.Item(43,5).Value = 0.9937578190994688
$ws.Cells.Item(44,1).Value = 27
$ws.Cells.Item(44,2).Value = 0.17338
$ws.Cells.Item(44,3).Value = 0.1657037
$ws.Cells.Item(44,4).Value = 0.007676299999999997
$ws.Cells.Item(44,5).Value = 0.9923821766970207
$ws.Cells.Item(45,1).Value = 37
$ws.Cells.Item(45,2).Value = 0.17029
$ws.Cells.Item(45,3).Value = 0.1616971
$ws.Cells.Item(45,4).Value = 0.008592900000000014
$ws.Cells.Item(45,5).Value = 0.991480308854048
$ws.Cells.Item(46,1).Value = 3
$ws.Cells.Item(46,2).Value = 0.11606
$ws.Cells.Item(46,3).Value = 0.1289021999999999
$ws.Cells.Item(46,4).Value = 0.01284219999999991
$ws.Cells.Item(46,5).Value = 0.9873206309926661
$ws.Cells.Item(47,1).Value = 31
$ws.Cells.Item(47,2).Value = 0.22902
$ws.Cells.Item(47,3).Value = 0.2136908
$ws.Cells.Item(47,4).Value = 0.01532919999999999
$ws.Cells.Item(47,5).Value = 0.9849022366341872
$ws.Cells.Item(48,1).Value = 45
$ws.Cells.Item(48,2).Value = 0.18125
$ws.Cells.Item(48,3).Value = 0.1618987
$ws.Cells.Item(48,4).Value = 0.01935129999999999
$ws.Cells.Item(48,5).Value = 0.9810160638437405
$ws.Cells.Item(49,1).Value = 2
$ws.Cells.Item(49,2).Value = 0.18856
$ws.Cells.Item(49,3).Value = 0.1682296
$ws.Cells.Item(49,4).Value = 0.02033040000000003
$ws.Cells.Item(49,5).Value = 0.9800746895319399
$ws.Cells.Item(50,1).Value = 34
$ws.Cells.Item(50,2).Value = 0.1377
$ws.Cells.Item(50,3).Value = 0.1646329666666667
$ws.Cells.Item(50,4).Value = 0.0269329666666667
$ws.Cells.Item(50,5).Value = 0.9737733936479918
$ws.Cells.Item(51,1).Value = 1
$ws.Cells.Item(51,2).Value = 0.22368
$ws.Cells.Item(51,3).Value = 0.1960796
$ws.Cells.Item(51,4).Value = 0.0276004
$ws.Cells.Item(51,5).Value = 0.9731409213153284
